# Update wing size (WSa/WSb) and body size (BSa/BSb) measurements for
# subjects 601-605 (rows 51-55), which previously had no data recorded
# in columns L:O.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newData = @{
    51 = @(15.92, 8.52,              15.89, 5.75)
    52 = @(15.75, 8.83,              15.3,  5.36)
    53 = @(15.68, 8.7100000000000009, 15.91, 5.19)
    54 = @(16.04, 8.75,              15.69, 5.61)
    55 = @(14.36, 8.4600000000000009, 15.35, 5.44)
}

foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    $ws.Range("L$row").Value = $vals[0]
    $ws.Range("M$row").Value = $vals[1]
    $ws.Range("N$row").Value = $vals[2]
    $ws.Range("O$row").Value = $vals[3]
}

# Reflect that the edit was made while rows 51:55 were selected (matches
# the saved sheet view of the authored workbook).
[void]$ws.Range("A51:XFD55").Select()
